$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$siteName = "Huxton"
$dateStr = "04.09.2014"
$runningFlag = "N"

# New data rows appended for the 04.09.2014 Huxton transect survey.
# Column B holds dates as literal text (matching the existing column's
# convention), so we briefly force a text format to stop Excel from
# auto-converting the date-like string into a date serial, then clear
# the format back to the default so the cell ends up styled like its
# neighbours.

$ws.Cells.Item(551, 1).Value2 = $siteName
$ws.Cells.Item(551, 2).NumberFormat = "@"
$ws.Cells.Item(551, 2).Value2 = $dateStr
$ws.Cells.Item(551, 2).ClearFormats()
$ws.Cells.Item(551, 3).Value2 = 1
$ws.Cells.Item(551, 4).Value2 = 0
$ws.Cells.Item(551, 5).Value2 = 0.86581018518518515
$ws.Cells.Item(551, 5).NumberFormat = "h:mm:ss"
$ws.Cells.Item(551, 6).Value2 = 10.5
$ws.Cells.Item(551, 7).Value2 = 0
$ws.Cells.Item(551, 8).Value2 = $runningFlag

$ws.Cells.Item(552, 1).Value2 = $siteName
$ws.Cells.Item(552, 2).NumberFormat = "@"
$ws.Cells.Item(552, 2).Value2 = $dateStr
$ws.Cells.Item(552, 2).ClearFormats()
$ws.Cells.Item(552, 3).Value2 = 1
$ws.Cells.Item(552, 4).Value2 = 100
$ws.Cells.Item(552, 5).Value2 = 0.86834490740740744
$ws.Cells.Item(552, 5).NumberFormat = "h:mm:ss"
$ws.Cells.Item(552, 6).Value2 = 10.5
$ws.Cells.Item(552, 7).Value2 = 0
$ws.Cells.Item(552, 8).Value2 = $runningFlag

$ws.Cells.Item(553, 1).Value2 = $siteName
$ws.Cells.Item(553, 2).NumberFormat = "@"
$ws.Cells.Item(553, 2).Value2 = $dateStr
$ws.Cells.Item(553, 2).ClearFormats()
$ws.Cells.Item(553, 3).Value2 = 1
$ws.Cells.Item(553, 4).Value2 = 200
$ws.Cells.Item(553, 5).Value2 = 0.86996527777777777
$ws.Cells.Item(553, 5).NumberFormat = "h:mm:ss"
$ws.Cells.Item(553, 6).Value2 = 10.5
$ws.Cells.Item(553, 7).Value2 = 0
$ws.Cells.Item(553, 8).Value2 = $runningFlag

$ws.Cells.Item(554, 1).Value2 = $siteName
$ws.Cells.Item(554, 2).NumberFormat = "@"
$ws.Cells.Item(554, 2).Value2 = $dateStr
$ws.Cells.Item(554, 2).ClearFormats()
$ws.Cells.Item(554, 3).Value2 = 1
$ws.Cells.Item(554, 4).Value2 = 300
$ws.Cells.Item(554, 5).Value2 = 0.87155092592592587
$ws.Cells.Item(554, 5).NumberFormat = "h:mm:ss"
$ws.Cells.Item(554, 6).Value2 = 10.5
$ws.Cells.Item(554, 7).Value2 = 0
$ws.Cells.Item(554, 8).Value2 = $runningFlag

$ws.Cells.Item(555, 1).Value2 = $siteName
$ws.Cells.Item(555, 2).NumberFormat = "@"
$ws.Cells.Item(555, 2).Value2 = $dateStr
$ws.Cells.Item(555, 2).ClearFormats()
$ws.Cells.Item(555, 3).Value2 = 1
$ws.Cells.Item(555, 4).Value2 = 400
$ws.Cells.Item(555, 5).Value2 = 0.87328703703703703
$ws.Cells.Item(555, 5).NumberFormat = "h:mm:ss"
$ws.Cells.Item(555, 6).Value2 = 10.5
$ws.Cells.Item(555, 7).Value2 = 0
$ws.Cells.Item(555, 8).Value2 = $runningFlag

$ws.Cells.Item(556, 1).Value2 = $siteName
$ws.Cells.Item(556, 2).NumberFormat = "@"
$ws.Cells.Item(556, 2).Value2 = $dateStr
$ws.Cells.Item(556, 2).ClearFormats()
$ws.Cells.Item(556, 3).Value2 = 1
$ws.Cells.Item(556, 4).Value2 = 500
$ws.Cells.Item(556, 5).Value2 = 0.87510416666666668
$ws.Cells.Item(556, 5).NumberFormat = "h:mm:ss"
$ws.Cells.Item(556, 6).Value2 = 10.5
$ws.Cells.Item(556, 7).Value2 = 0
$ws.Cells.Item(556, 8).Value2 = $runningFlag

$ws.Cells.Item(557, 1).Value2 = $siteName
$ws.Cells.Item(557, 2).NumberFormat = "@"
$ws.Cells.Item(557, 2).Value2 = $dateStr
$ws.Cells.Item(557, 2).ClearFormats()
$ws.Cells.Item(557, 3).Value2 = 1
$ws.Cells.Item(557, 4).Value2 = 600
$ws.Cells.Item(557, 5).Value2 = 0.87699074074074079
$ws.Cells.Item(557, 5).NumberFormat = "h:mm:ss"
$ws.Cells.Item(557, 6).Value2 = 10.5
$ws.Cells.Item(557, 7).Value2 = 0
$ws.Cells.Item(557, 8).Value2 = $runningFlag

$ws.Cells.Item(558, 1).Value2 = $siteName
$ws.Cells.Item(558, 2).NumberFormat = "@"
$ws.Cells.Item(558, 2).Value2 = $dateStr
$ws.Cells.Item(558, 2).ClearFormats()
$ws.Cells.Item(558, 3).Value2 = 2
$ws.Cells.Item(558, 4).Value2 = 0
$ws.Cells.Item(558, 5).Value2 = 0.91931712962962964
$ws.Cells.Item(558, 5).NumberFormat = "h:mm:ss"
$ws.Cells.Item(558, 6).Value2 = 10.5
$ws.Cells.Item(558, 7).Value2 = 0
$ws.Cells.Item(558, 8).Value2 = $runningFlag

$ws.Cells.Item(559, 1).Value2 = $siteName
$ws.Cells.Item(559, 2).NumberFormat = "@"
$ws.Cells.Item(559, 2).Value2 = $dateStr
$ws.Cells.Item(559, 2).ClearFormats()
$ws.Cells.Item(559, 3).Value2 = 2
$ws.Cells.Item(559, 4).Value2 = 100
$ws.Cells.Item(559, 5).Value2 = 0.92143518518518519
$ws.Cells.Item(559, 5).NumberFormat = "h:mm:ss"
$ws.Cells.Item(559, 6).Value2 = 10.5
$ws.Cells.Item(559, 7).Value2 = 0
$ws.Cells.Item(559, 8).Value2 = $runningFlag

$ws.Cells.Item(560, 1).Value2 = $siteName
$ws.Cells.Item(560, 2).NumberFormat = "@"
$ws.Cells.Item(560, 2).Value2 = $dateStr
$ws.Cells.Item(560, 2).ClearFormats()
$ws.Cells.Item(560, 3).Value2 = 2
$ws.Cells.Item(560, 4).Value2 = 200
$ws.Cells.Item(560, 5).Value2 = 0.92358796296296297
$ws.Cells.Item(560, 5).NumberFormat = "h:mm:ss"
$ws.Cells.Item(560, 6).Value2 = 10.5
$ws.Cells.Item(560, 7).Value2 = 0
$ws.Cells.Item(560, 8).Value2 = $runningFlag

$ws.Cells.Item(561, 1).Value2 = $siteName
$ws.Cells.Item(561, 2).NumberFormat = "@"
$ws.Cells.Item(561, 2).Value2 = $dateStr
$ws.Cells.Item(561, 2).ClearFormats()
$ws.Cells.Item(561, 3).Value2 = 2
$ws.Cells.Item(561, 4).Value2 = 300
$ws.Cells.Item(561, 5).Value2 = 0.92497685185185186
$ws.Cells.Item(561, 5).NumberFormat = "h:mm:ss"
$ws.Cells.Item(561, 6).Value2 = 10.5
$ws.Cells.Item(561, 7).Value2 = 0
$ws.Cells.Item(561, 8).Value2 = $runningFlag

$ws.Cells.Item(562, 1).Value2 = $siteName
$ws.Cells.Item(562, 2).NumberFormat = "@"
$ws.Cells.Item(562, 2).Value2 = $dateStr
$ws.Cells.Item(562, 2).ClearFormats()
$ws.Cells.Item(562, 3).Value2 = 2
$ws.Cells.Item(562, 4).Value2 = 400
$ws.Cells.Item(562, 5).Value2 = 0.9264930555555555
$ws.Cells.Item(562, 5).NumberFormat = "h:mm:ss"
$ws.Cells.Item(562, 6).Value2 = 10.5
$ws.Cells.Item(562, 7).Value2 = 0
$ws.Cells.Item(562, 8).Value2 = $runningFlag

$ws.Cells.Item(563, 1).Value2 = $siteName
$ws.Cells.Item(563, 2).NumberFormat = "@"
$ws.Cells.Item(563, 2).Value2 = $dateStr
$ws.Cells.Item(563, 2).ClearFormats()
$ws.Cells.Item(563, 3).Value2 = 2
$ws.Cells.Item(563, 4).Value2 = 500
$ws.Cells.Item(563, 5).Value2 = 0.92791666666666661
$ws.Cells.Item(563, 5).NumberFormat = "h:mm:ss"
$ws.Cells.Item(563, 6).Value2 = 10.5
$ws.Cells.Item(563, 7).Value2 = 0
$ws.Cells.Item(563, 8).Value2 = $runningFlag

$ws.Cells.Item(564, 1).Value2 = $siteName
$ws.Cells.Item(564, 2).NumberFormat = "@"
$ws.Cells.Item(564, 2).Value2 = $dateStr
$ws.Cells.Item(564, 2).ClearFormats()
$ws.Cells.Item(564, 3).Value2 = 2
$ws.Cells.Item(564, 4).Value2 = 600
$ws.Cells.Item(564, 5).Value2 = 0.92962962962962958
$ws.Cells.Item(564, 5).NumberFormat = "h:mm:ss"
$ws.Cells.Item(564, 6).Value2 = 10.5
$ws.Cells.Item(564, 7).Value2 = 0
$ws.Cells.Item(564, 8).Value2 = $runningFlag

$ws.Cells.Item(565, 1).Value2 = $siteName
$ws.Cells.Item(565, 2).NumberFormat = "@"
$ws.Cells.Item(565, 2).Value2 = $dateStr
$ws.Cells.Item(565, 2).ClearFormats()
$ws.Cells.Item(565, 3).Value2 = 3
$ws.Cells.Item(565, 4).Value2 = 0
$ws.Cells.Item(565, 5).Value2 = 0.90055555555555555
$ws.Cells.Item(565, 5).NumberFormat = "h:mm:ss"
$ws.Cells.Item(565, 6).Value2 = 10.5
$ws.Cells.Item(565, 7).Value2 = 0
$ws.Cells.Item(565, 8).Value2 = $runningFlag

$ws.Cells.Item(566, 1).Value2 = $siteName
$ws.Cells.Item(566, 2).NumberFormat = "@"
$ws.Cells.Item(566, 2).Value2 = $dateStr
$ws.Cells.Item(566, 2).ClearFormats()
$ws.Cells.Item(566, 3).Value2 = 3
$ws.Cells.Item(566, 4).Value2 = 100
$ws.Cells.Item(566, 5).Value2 = 0.90312500000000007
$ws.Cells.Item(566, 5).NumberFormat = "h:mm:ss"
$ws.Cells.Item(566, 6).Value2 = 10.5
$ws.Cells.Item(566, 7).Value2 = 0
$ws.Cells.Item(566, 8).Value2 = $runningFlag

$ws.Cells.Item(567, 1).Value2 = $siteName
$ws.Cells.Item(567, 2).NumberFormat = "@"
$ws.Cells.Item(567, 2).Value2 = $dateStr
$ws.Cells.Item(567, 2).ClearFormats()
$ws.Cells.Item(567, 3).Value2 = 3
$ws.Cells.Item(567, 4).Value2 = 200
$ws.Cells.Item(567, 5).Value2 = 0.90494212962962972
$ws.Cells.Item(567, 5).NumberFormat = "h:mm:ss"
$ws.Cells.Item(567, 6).Value2 = 10.5
$ws.Cells.Item(567, 7).Value2 = 0
$ws.Cells.Item(567, 8).Value2 = $runningFlag

$ws.Cells.Item(568, 1).Value2 = $siteName
$ws.Cells.Item(568, 2).NumberFormat = "@"
$ws.Cells.Item(568, 2).Value2 = $dateStr
$ws.Cells.Item(568, 2).ClearFormats()
$ws.Cells.Item(568, 3).Value2 = 3
$ws.Cells.Item(568, 4).Value2 = 300
$ws.Cells.Item(568, 5).Value2 = 0.90678240740740745
$ws.Cells.Item(568, 5).NumberFormat = "h:mm:ss"
$ws.Cells.Item(568, 6).Value2 = 10.5
$ws.Cells.Item(568, 7).Value2 = 0
$ws.Cells.Item(568, 8).Value2 = $runningFlag

$ws.Cells.Item(569, 1).Value2 = $siteName
$ws.Cells.Item(569, 2).NumberFormat = "@"
$ws.Cells.Item(569, 2).Value2 = $dateStr
$ws.Cells.Item(569, 2).ClearFormats()
$ws.Cells.Item(569, 3).Value2 = 3
$ws.Cells.Item(569, 4).Value2 = 400
$ws.Cells.Item(569, 5).Value2 = 0.90861111111111104
$ws.Cells.Item(569, 5).NumberFormat = "h:mm:ss"
$ws.Cells.Item(569, 6).Value2 = 10.5
$ws.Cells.Item(569, 7).Value2 = 0
$ws.Cells.Item(569, 8).Value2 = $runningFlag

$ws.Cells.Item(570, 1).Value2 = $siteName
$ws.Cells.Item(570, 2).NumberFormat = "@"
$ws.Cells.Item(570, 2).Value2 = $dateStr
$ws.Cells.Item(570, 2).ClearFormats()
$ws.Cells.Item(570, 3).Value2 = 3
$ws.Cells.Item(570, 4).Value2 = 500
$ws.Cells.Item(570, 5).Value2 = 0.91038194444444442
$ws.Cells.Item(570, 5).NumberFormat = "h:mm:ss"
$ws.Cells.Item(570, 6).Value2 = 10.5
$ws.Cells.Item(570, 7).Value2 = 0
$ws.Cells.Item(570, 8).Value2 = $runningFlag

$ws.Cells.Item(571, 1).Value2 = $siteName
$ws.Cells.Item(571, 2).NumberFormat = "@"
$ws.Cells.Item(571, 2).Value2 = $dateStr
$ws.Cells.Item(571, 2).ClearFormats()
$ws.Cells.Item(571, 3).Value2 = 3
$ws.Cells.Item(571, 4).Value2 = 600
$ws.Cells.Item(571, 5).Value2 = 0.91200231481481486
$ws.Cells.Item(571, 5).NumberFormat = "h:mm:ss"
$ws.Cells.Item(571, 6).Value2 = 10.5
$ws.Cells.Item(571, 7).Value2 = 0
$ws.Cells.Item(571, 8).Value2 = $runningFlag

$ws.Cells.Item(572, 1).Value2 = $siteName
$ws.Cells.Item(572, 2).NumberFormat = "@"
$ws.Cells.Item(572, 2).Value2 = $dateStr
$ws.Cells.Item(572, 2).ClearFormats()
$ws.Cells.Item(572, 3).Value2 = 3
$ws.Cells.Item(572, 4).Value2 = 700
$ws.Cells.Item(572, 5).Value2 = 0.91361111111111104
$ws.Cells.Item(572, 5).NumberFormat = "h:mm:ss"
$ws.Cells.Item(572, 6).Value2 = 10.5
$ws.Cells.Item(572, 7).Value2 = 0
$ws.Cells.Item(572, 8).Value2 = $runningFlag

$ws.Cells.Item(573, 1).Value2 = $siteName
$ws.Cells.Item(573, 2).NumberFormat = "@"
$ws.Cells.Item(573, 2).Value2 = $dateStr
$ws.Cells.Item(573, 2).ClearFormats()
$ws.Cells.Item(573, 3).Value2 = 4
$ws.Cells.Item(573, 4).Value2 = 0
$ws.Cells.Item(573, 5).Value2 = 0.88351851851851848
$ws.Cells.Item(573, 5).NumberFormat = "h:mm:ss"
$ws.Cells.Item(573, 6).Value2 = 10.5
$ws.Cells.Item(573, 7).Value2 = 0
$ws.Cells.Item(573, 8).Value2 = $runningFlag

$ws.Cells.Item(574, 1).Value2 = $siteName
$ws.Cells.Item(574, 2).NumberFormat = "@"
$ws.Cells.Item(574, 2).Value2 = $dateStr
$ws.Cells.Item(574, 2).ClearFormats()
$ws.Cells.Item(574, 3).Value2 = 4
$ws.Cells.Item(574, 4).Value2 = 100
$ws.Cells.Item(574, 5).Value2 = 0.88615740740740734
$ws.Cells.Item(574, 5).NumberFormat = "h:mm:ss"
$ws.Cells.Item(574, 6).Value2 = 10.5
$ws.Cells.Item(574, 7).Value2 = 0
$ws.Cells.Item(574, 8).Value2 = $runningFlag

$ws.Cells.Item(575, 1).Value2 = $siteName
$ws.Cells.Item(575, 2).NumberFormat = "@"
$ws.Cells.Item(575, 2).Value2 = $dateStr
$ws.Cells.Item(575, 2).ClearFormats()
$ws.Cells.Item(575, 3).Value2 = 4
$ws.Cells.Item(575, 4).Value2 = 200
$ws.Cells.Item(575, 5).Value2 = 0.88800925925925922
$ws.Cells.Item(575, 5).NumberFormat = "h:mm:ss"
$ws.Cells.Item(575, 6).Value2 = 10.5
$ws.Cells.Item(575, 7).Value2 = 0
$ws.Cells.Item(575, 8).Value2 = $runningFlag

$ws.Cells.Item(576, 1).Value2 = $siteName
$ws.Cells.Item(576, 2).NumberFormat = "@"
$ws.Cells.Item(576, 2).Value2 = $dateStr
$ws.Cells.Item(576, 2).ClearFormats()
$ws.Cells.Item(576, 3).Value2 = 4
$ws.Cells.Item(576, 4).Value2 = 300
$ws.Cells.Item(576, 5).Value2 = 0.88974537037037038
$ws.Cells.Item(576, 5).NumberFormat = "h:mm:ss"
$ws.Cells.Item(576, 6).Value2 = 10.5
$ws.Cells.Item(576, 7).Value2 = 0
$ws.Cells.Item(576, 8).Value2 = $runningFlag

$ws.Cells.Item(577, 1).Value2 = $siteName
$ws.Cells.Item(577, 2).NumberFormat = "@"
$ws.Cells.Item(577, 2).Value2 = $dateStr
$ws.Cells.Item(577, 2).ClearFormats()
$ws.Cells.Item(577, 3).Value2 = 4
$ws.Cells.Item(577, 4).Value2 = 400
$ws.Cells.Item(577, 5).Value2 = 0.89134259259259263
$ws.Cells.Item(577, 5).NumberFormat = "h:mm:ss"
$ws.Cells.Item(577, 6).Value2 = 10.5
$ws.Cells.Item(577, 7).Value2 = 0
$ws.Cells.Item(577, 8).Value2 = $runningFlag

$ws.Cells.Item(578, 1).Value2 = $siteName
$ws.Cells.Item(578, 2).NumberFormat = "@"
$ws.Cells.Item(578, 2).Value2 = $dateStr
$ws.Cells.Item(578, 2).ClearFormats()
$ws.Cells.Item(578, 3).Value2 = 4
$ws.Cells.Item(578, 4).Value2 = 500
$ws.Cells.Item(578, 5).Value2 = 0.89307870370370368
$ws.Cells.Item(578, 5).NumberFormat = "h:mm:ss"
$ws.Cells.Item(578, 6).Value2 = 10.5
$ws.Cells.Item(578, 7).Value2 = 0
$ws.Cells.Item(578, 8).Value2 = $runningFlag

$ws.Cells.Item(579, 1).Value2 = $siteName
$ws.Cells.Item(579, 2).NumberFormat = "@"
$ws.Cells.Item(579, 2).Value2 = $dateStr
$ws.Cells.Item(579, 2).ClearFormats()
$ws.Cells.Item(579, 3).Value2 = 4
$ws.Cells.Item(579, 4).Value2 = 600
$ws.Cells.Item(579, 5).Value2 = 0.89482638888888888
$ws.Cells.Item(579, 5).NumberFormat = "h:mm:ss"
$ws.Cells.Item(579, 6).Value2 = 10.5
$ws.Cells.Item(579, 7).Value2 = 0
$ws.Cells.Item(579, 8).Value2 = $runningFlag

# Move the selection to reflect where the user ended up after entering
# the new rows (the header row stays frozen via the pre-existing pane).
[void]$ws.Activate()
[void]$ws.Range("E580").Select()

Write-Output "Added rows 551-579 (Huxton, 04.09.2014)."
